$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "preview" row driving the formulas in B4:B6 (and mirrored later into A13:A15)
$ws.Range("A2").Value = "INTNL_ISSS"
$ws.Range("B2").Value = "International Student and Scholar Services"
$ws.Range("C2").Value = "INTNLdata()"
$ws.Range("D2").Value = "International Student and Scholar Services"

# Row 10: full resource record, replacing the old "University Interfaith Center" entry
$ws.Range("A10").Value = "International Student and Scholar Services"
$ws.Range("B10").Value = "Science Library G-40  "
$ws.Range("C10").Value = "tel:(518)5918172"
$ws.Range("D10").Value = "https://www.albany.edu/isss/"
$ws.Range("E10").Value = "M, Tues, Th, F: 1 PM - 3:30 PM"
$ws.Range("F10").Value = "ISSS@albany.edu"
$ws.Range("G10").Value = "INTNL"
$ws.Range("H10").Value = "On"
$ws.Range("I10").Value = "Advising - workshops"

# Remove hyperlinks that used to sit on C10/D10/F10 for the old record, and
# drop the underlined-blue hyperlink look in favour of the plain row font
$ws.Hyperlinks.Delete()
$ws.Range("C10").Font.Name = "Times New Roman"
$ws.Range("C10").Font.Size = 12
$ws.Range("C10").Font.Underline = $false
$ws.Range("D10").Font.Name = "Times New Roman"
$ws.Range("D10").Font.Size = 12
$ws.Range("D10").Font.Underline = $false
$ws.Range("F10").Font.Name = "Times New Roman"
$ws.Range("F10").Font.Size = 12
$ws.Range("F10").Font.Underline = $false

# Paste the recalculated formula results (B4:B6) as literal values into A13:A15
$ws.Range("A13").Value = $ws.Range("B4").Value2
$ws.Range("A14").Value = $ws.Range("B5").Value2
$ws.Range("A15").Value = $ws.Range("B6").Value2

$ws.Range("A13:A15").Select()
